$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 12955
$ws.Range("F3").Value = 7255
$ws.Range("F4").Value = 267
$ws.Range("F5").Value = 12
$ws.Range("F7").Value = 1024
$ws.Range("G7").Value = "不可售"
$ws.Range("F10").Value = 1048
$ws.Range("F12").Value = 70
$ws.Range("F13").Value = 1035
$ws.Range("F14").Value = 515
$ws.Range("F15").Value = 274
$ws.Range("F16").Value = 382
$ws.Range("F19").Value = 323
$ws.Range("F21").Value = 257
$ws.Range("F22").Value = 402
$ws.Range("F23").Value = 5307
$ws.Range("F25").Value = 1475
$ws.Range("F26").Value = 328
$ws.Range("F27").Value = 1804
$ws.Range("F28").Value = 1804
$ws.Range("F29").Value = 104
$ws.Range("F31").Value = 1398
$ws.Range("F34").Value = 613
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 78
$ws.Range("F12").Value = 13
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9328
$ws.Range("F4").Value = 2069
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9328
$ws.Range("F4").Value = 2069
$ws.Range("F5").Value = 12955
$ws.Range("F6").Value = 7255
$ws.Range("F8").Value = 1024
$ws.Range("G8").Value = "不可售"
$ws.Range("F11").Value = 1048
$ws.Range("F13").Value = 70
$ws.Range("F14").Value = 1035
$ws.Range("F15").Value = 515
$ws.Range("F16").Value = 274
$ws.Range("F17").Value = 382
$ws.Range("F20").Value = 323
$ws.Range("F25").Value = 257
$ws.Range("F26").Value = 402
$ws.Range("F27").Value = 5307
$ws.Range("F29").Value = 1475
$ws.Range("F32").Value = 328
$ws.Range("F34").Value = 1804
$ws.Range("F35").Value = 1804
$ws.Range("F36").Value = 104
$ws.Range("F38").Value = 1398
$ws.Range("F41").Value = 613
$ws.Range("F42").Value = 13
